$d = $word.ActiveDocument

# Replace the month abbreviation "Dec" -> "Jan" in the title line.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Dec", $true, $true, $false, $false, $false, $true, 1, $false, "Jan", 2)

# Replace the year "2023" -> "2024" in the title line.
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("2023", $true, $true, $false, $false, $false, $true, 1, $false, "2024", 2)
